$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 05:51:32"
$ws1.Range("A3").Value = "Total filas: 43"

# New data block for rows 31-48 (replaces old rows 31-36, extends to 48)
$sheet1Data = @(
    @("05:51:32", "06:59", "14_ABASTO", 68, "LP1912"),
    @("05:20:00", "07:00", "10_OLMOS", 100, "LP1912"),
    @("05:20:00", "07:00", "14_ABASTO", 100, "LP1912"),
    @("05:51:32", "07:04", "23_HERNANDEZ", 73, "LP1912"),
    @("05:20:00", "07:05", "15_ABASTO", 105, "LP1912"),
    @("05:20:00", "07:07", "225_GOMEZ", 107, "LP1912"),
    @("05:51:32", "07:11", "215A_EL PATO", 80, "LP1912"),
    @("05:20:00", "07:12", "215A_EL PATO", 112, "LP1912"),
    @("05:51:32", "07:15", "11_ETCHEVERRY", 84, "LP1912"),
    @("05:20:00", "07:16", "11_ETCHEVERRY", 116, "LP1912"),
    @("05:51:32", "07:21", "26_HERNANDEZ", 90, "LP1912"),
    @("05:51:32", "07:28", "10_OLMOS", 97, "LP1912"),
    @("05:51:32", "07:31", "11_ETCHEVERRY", 100, "LP1912"),
    @("05:51:32", "07:31", "16_SANTA ANA", 100, "LP1912"),
    @("05:51:32", "07:32", "84_COLONIA URQUIZA-ESC 49", 101, "LP1912"),
    @("05:51:32", "07:36", "27_EL RETIRO", 105, "LP1912"),
    @("05:51:32", "07:39", "10_OLMOS", 108, "LP1912"),
    @("05:51:32", "07:47", "14_ABASTO", 116, "LP1912")
)

$r = 31
foreach ($row in $sheet1Data) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 05:51:32"
$ws2.Range("A3").Value = "Total filas: 10"

# Shift old row 14 down to row 15, then insert new row 14
$ws2.Cells.Item(15, 1).Value = "05:20:00"
$ws2.Cells.Item(15, 2).Value = "07:12"
$ws2.Cells.Item(15, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(15, 4).Value = 112
$ws2.Cells.Item(15, 5).Value = "LP1912"

$ws2.Cells.Item(14, 1).Value = "05:51:32"
$ws2.Cells.Item(14, 2).Value = "07:11"
$ws2.Cells.Item(14, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(14, 4).Value = 80
$ws2.Cells.Item(14, 5).Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 05:51:32"
$ws3.Range("A3").Value = "Total filas: 8"

# Append new row 13
$ws3.Cells.Item(13, 1).Value = "05:51:32"
$ws3.Cells.Item(13, 2).Value = "07:35"
$ws3.Cells.Item(13, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(13, 4).Value = 104
$ws3.Cells.Item(13, 5).Value = "L6173"
